# Apply cryptos list price/volume updates (commit: "Updated cryptos list on Mon Jun  3 19:49:07 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.079.89'
$ws.Range("E2").Value = '  +1.98%  '

# Row 3
$ws.Range("D3").Value = '3.772.01'
$ws.Range("E3").Value = '  -0.15%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''625.32'  # force text (numeric-looking string)
$ws.Range("E5").Value = '  +4.38%  '

# Row 6
$ws.Range("D6").Value = '''166.60'  # force text (numeric-looking string)
$ws.Range("E6").Value = '  +2.24%  '

# Row 7
$ws.Range("D7").Value = '3.771.07'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("E8").Value = '  -0.11%  '

# Row 9
$ws.Range("E9").Value = '  +1.61%  '

# Row 10
$ws.Range("E10").Value = '  +2.78%  '

# Row 11
$ws.Range("E11").Value = '  +3.02%  '

# Row 12
$ws.Range("D12").Value = '''6.70'  # force text (numeric-looking string)
$ws.Range("E12").Value = '  +1.13%  '

# Row 13
$ws.Range("E13").Value = '  +1.34%  '

# Row 14
$ws.Range("E14").Value = '  +1.87%  '

# Row 15
$ws.Range("D15").Value = '4.410.96'
$ws.Range("E15").Value = '  +0.03%  '

# Row 16
$ws.Range("D16").Value = '3.774.35'
$ws.Range("E16").Value = '  +0.20%  '

# Row 17
$ws.Range("D17").Value = '69.101.71'
$ws.Range("E17").Value = '  +2.06%  '

# Row 18
$ws.Range("D18").Value = '''17.65'  # force text (numeric-looking string)
$ws.Range("E18").Value = '  -2.74%  '

# Row 19
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '''0.114'  # force text (numeric-looking string)
$ws.Range("E19").Value = '  -0.89%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''7.09'  # force text (numeric-looking string)
$ws.Range("E20").Value = '  +1.21%  '

# Row 21
$ws.Range("D21").Value = '''468.37'  # force text (numeric-looking string)
$ws.Range("E21").Value = '  +2.55%  '

# Row 22
$ws.Range("E22").Value = '  +1.60%  '

# Row 23
$ws.Range("E23").Value = '  +2.44%  '

# Row 24
$ws.Range("E24").Value = '  +4.69%  '

# Row 25
$ws.Range("D25").Value = '''83.11'  # force text (numeric-looking string)
$ws.Range("E25").Value = '  +0.20%  '

# Row 26
$ws.Range("D26").Value = '''12.08'  # force text (numeric-looking string)
$ws.Range("E26").Value = '  +2.21%  '

# Row 27
$ws.Range("E27").Value = '  +3.90%  '

# Row 28
$ws.Range("D28").Value = '''10.04'  # force text (numeric-looking string)
$ws.Range("E28").Value = '  +2.08%  '

# Row 29
$ws.Range("E29").Value = '  -0.06%  '

# Row 30
$ws.Range("D30").Value = '3.922.78'
$ws.Range("E30").Value = '  +0.11%  '

# Row 32
$ws.Range("E32").Value = '  +2.39%  '

# Row 33
$ws.Range("E33").Value = '  +0.75%  '

# Row 34
$ws.Range("D34").Value = '''28.76'  # force text (numeric-looking string)
$ws.Range("E34").Value = '  -0.33%  '

# Row 35
$ws.Range("E35").Value = '  +0.07%  '

# Row 36
$ws.Range("D36").Value = '''8.99'  # force text (numeric-looking string)
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("D37").Value = '3.724.73'
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("D38").Value = '''0.163'  # force text (numeric-looking string)
$ws.Range("E38").Value = '  +13.49%  '

# Row 39
$ws.Range("E39").Value = '  +2.74%  '

# Row 40
$ws.Range("D40").Value = '''3.45'  # force text (numeric-looking string)
$ws.Range("E40").Value = '  +8.90%  '

# Row 41
$ws.Range("E41").Value = '  +1.04%  '

# Row 42
$ws.Range("D42").Value = '''0.967'  # force text (numeric-looking string)
$ws.Range("E42").Value = '  -1.10%  '

# Row 43
$ws.Range("D43").Value = '''1.00'  # force text (numeric-looking string)

# Row 44
$ws.Range("E44").Value = '  -0.04%  '

# Row 45
$ws.Range("B45").Value = 'Arweave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D45").Value = '''43.27'  # force text (numeric-looking string)
$ws.Range("E45").Value = '  -0.35%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '''0.298'  # force text (numeric-looking string)
$ws.Range("E46").Value = '  +1.21%  '

# Row 47
$ws.Range("D47").Value = '''152.42'  # force text (numeric-looking string)
$ws.Range("E47").Value = '  +0.49%  '

# Row 48
$ws.Range("E48").Value = '  +4.43%  '

# Row 49
$ws.Range("D49").Value = '''46.64'  # force text (numeric-looking string)
$ws.Range("E49").Value = '  -1.19%  '

# Row 50
$ws.Range("E50").Value = '  +1.80%  '

# Row 51
$ws.Range("E51").Value = '  +0.02%  '
